$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to plain text so numeric-looking strings
# (e.g. "214.86") are not reinterpreted as numbers/dates by Excel,
# then restore the default "Normal" style so no stray formatting is left behind.
function Set-TextValue {
    param($cell, [string]$text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "26.986.75"
$ws.Range("D3").Value = "1.671.77"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").Value = "  +0.10%  "
Set-TextValue $ws.Range("D5") "214.86"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("E6").Value = "  -1.17%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +1.76%  "
Set-TextValue $ws.Range("D9") "21.40"
$ws.Range("E9").Value = "  +5.45%  "
Set-TextValue $ws.Range("D10") "0.0621"
$ws.Range("E10").Value = "  -0.06%  "
Set-TextValue $ws.Range("D11") "0.0887"
$ws.Range("E11").Value = "  -0.56%  "
$ws.Range("D12").Value = "1.907.61"
$ws.Range("E12").Value = "  +0.20%  "
$ws.Range("D13").Value = "1.664.18"
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("E14").Value = "  +0.69%  "
$ws.Range("E15").Value = "  +1.42%  "
Set-TextValue $ws.Range("D16") "66.08"
$ws.Range("E16").Value = "  +0.66%  "
$ws.Range("D17").Value = "26.988.12"
$ws.Range("E17").Value = "  +0.26%  "
Set-TextValue $ws.Range("D18") "8.18"
$ws.Range("E18").Value = "  +2.86%  "
Set-TextValue $ws.Range("D19") "234.45"
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("D20").Value = "0.0₃0734"
$ws.Range("E20").Value = "  +0.30%  "
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("E22").Value = "  +1.35%  "
$ws.Range("E23").Value = "  +0.79%  "
$ws.Range("E24").Value = "  -2.49%  "
Set-TextValue $ws.Range("D25") "146.56"
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("E26").Value = "  +1.91%  "
$ws.Range("E27").Value = "  +3.02%  "
$ws.Range("E28").Value = "  -0.55%  "
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("E30").Value = "  +0.56%  "
$ws.Range("E31").Value = "  -0.21%  "
$ws.Range("E32").Value = "  +0.51%  "
$ws.Range("D33").Value = "1.533.94"
$ws.Range("E33").Value = "  +5.71%  "
$ws.Range("E34").Value = "  +0.59%  "
$ws.Range("E35").Value = "  +3.62%  "
$ws.Range("E36").Value = "  -1.15%  "
$ws.Range("E37").Value = "  +1.29%  "
$ws.Range("E38").Value = "  +2.17%  "
Set-TextValue $ws.Range("D39") "0.908"
$ws.Range("E39").Value = "  +0.46%  "
$ws.Range("E40").Value = "  +5.80%  "
$ws.Range("E41").Value = "  +0.01%  "
Set-TextValue $ws.Range("D42") "67.57"
$ws.Range("E42").Value = "  +2.14%  "
$ws.Range("E43").Value = "  -3.51%  "
$ws.Range("E44").Value = "  -2.38%  "
$ws.Range("D45").Value = "1.816.18"
$ws.Range("E45").Value = "  +0.59%  "
$ws.Range("E46").Value = "  -0.12%  "
Set-TextValue $ws.Range("D47") "90.33"
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D48") "1.54"
$ws.Range("E48").Value = "  -0.19%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D49") "0.103"
$ws.Range("E49").Value = "  +1.75%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D50") "7.99"
$ws.Range("E50").Value = "  +5.79%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D51") "0.0509"
$ws.Range("E51").Value = "  +0.17%  "

Write-Host "Applied cryptos list update ($([int]79) cells)"
